$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi LR-pairs (Col9a3 -> Mag) results, recomputed following Dr Hou's advice.
# Rows 2-4 are refreshed with new values and rows 5-9 are added so that every
# Sending cluster (ECs, FAPs, M2, sCs) x Target cluster (M2, sCs) combination
# for the Col9a3 (ligand) -> Mag (receptor) pair is represented.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col9a3"
$ws.Range("C2").Value = "Mag"
$ws.Range("D2").Value = "M2"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6610796666666666
$ws.Range("H2").Value = 1.983239
$ws.Range("I2").Value = 0.2677179076566869
$ws.Range("J2").Value = 0.2677179076566869
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7885686666666668
$ws.Range("N2").Value = 2.365706
$ws.Range("O2").Value = 0.4566863346753138
$ws.Range("P2").Value = 0.4566863346753137
$ws.Range("Q2").Value = 0.5213067113037778
$ws.Range("R2").Value = 4.691760401734
$ws.Range("S2").Value = 0.1222631099746765
$ws.Range("T2").Value = 0.1222631099746764

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col9a3"
$ws.Range("C3").Value = "Mag"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6610796666666666
$ws.Range("H3").Value = 1.983239
$ws.Range("I3").Value = 0.2677179076566869
$ws.Range("J3").Value = 0.2677179076566869
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.9381496666666668
$ws.Range("N3").Value = 2.814449
$ws.Range("O3").Value = 0.5433136653246862
$ws.Range("P3").Value = 0.5433136653246862
$ws.Range("Q3").Value = 0.6201916689234445
$ws.Range("R3").Value = 5.581725020311
$ws.Range("S3").Value = 0.1454547976820104
$ws.Range("T3").Value = 0.1454547976820104

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Col9a3"
$ws.Range("C4").Value = "Mag"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.310047666666666
$ws.Range("H4").Value = 3.930143
$ws.Range("I4").Value = 0.5305309449600247
$ws.Range("J4").Value = 0.5305309449600246
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7885686666666668
$ws.Range("N4").Value = 2.365706
$ws.Range("O4").Value = 0.4566863346753138
$ws.Range("P4").Value = 0.4566863346753137
$ws.Range("Q4").Value = 1.033062541773111
$ws.Range("R4").Value = 9.297562875958
$ws.Range("S4").Value = 0.2422862326856243
$ws.Range("T4").Value = 0.2422862326856242

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col9a3"
$ws.Range("C5").Value = "Mag"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.310047666666666
$ws.Range("H5").Value = 3.930143
$ws.Range("I5").Value = 0.5305309449600247
$ws.Range("J5").Value = 0.5305309449600246
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.9381496666666668
$ws.Range("N5").Value = 2.814449
$ws.Range("O5").Value = 0.5433136653246862
$ws.Range("P5").Value = 0.5433136653246862
$ws.Range("Q5").Value = 1.229020781800778
$ws.Range("R5").Value = 11.061187036207
$ws.Range("S5").Value = 0.2882447122744003
$ws.Range("T5").Value = 0.2882447122744003

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Col9a3"
$ws.Range("C6").Value = "Mag"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.311447
$ws.Range("H6").Value = 0.934341
$ws.Range("I6").Value = 0.1261269153933825
$ws.Range("J6").Value = 0.1261269153933825
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7885686666666668
$ws.Range("N6").Value = 2.365706
$ws.Range("O6").Value = 0.4566863346753138
$ws.Range("P6").Value = 0.4566863346753137
$ws.Range("Q6").Value = 0.2455973455273334
$ws.Range("R6").Value = 2.210376109746
$ws.Range("S6").Value = 0.05760043869490727
$ws.Range("T6").Value = 0.05760043869490726

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Col9a3"
$ws.Range("C7").Value = "Mag"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.311447
$ws.Range("H7").Value = 0.934341
$ws.Range("I7").Value = 0.1261269153933825
$ws.Range("J7").Value = 0.1261269153933825
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.9381496666666668
$ws.Range("N7").Value = 2.814449
$ws.Range("O7").Value = 0.5433136653246862
$ws.Range("P7").Value = 0.5433136653246862
$ws.Range("Q7").Value = 0.2921838992343334
$ws.Range("R7").Value = 2.629655093109
$ws.Range("S7").Value = 0.06852647669847522
$ws.Range("T7").Value = 0.06852647669847522

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col9a3"
$ws.Range("C8").Value = "Mag"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.18674
$ws.Range("H8").Value = 0.5602199999999999
$ws.Range("I8").Value = 0.07562423198990599
$ws.Range("J8").Value = 0.07562423198990596
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.7885686666666668
$ws.Range("N8").Value = 2.365706
$ws.Range("O8").Value = 0.4566863346753138
$ws.Range("P8").Value = 0.4566863346753137
$ws.Range("Q8").Value = 0.1472573128133333
$ws.Range("R8").Value = 1.32531581532
$ws.Range("S8").Value = 0.03453655332010578
$ws.Range("T8").Value = 0.03453655332010576

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col9a3"
$ws.Range("C9").Value = "Mag"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.18674
$ws.Range("H9").Value = 0.5602199999999999
$ws.Range("I9").Value = 0.07562423198990599
$ws.Range("J9").Value = 0.07562423198990596
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.9381496666666668
$ws.Range("N9").Value = 2.814449
$ws.Range("O9").Value = 0.5433136653246862
$ws.Range("P9").Value = 0.5433136653246862
$ws.Range("Q9").Value = 0.1751900687533333
$ws.Range("R9").Value = 1.57671061878
$ws.Range("S9").Value = 0.04108767866980021
$ws.Range("T9").Value = 0.0410876786698002
